$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing C column dates (rows 4-15) to the new series -----------
$cVals = @{
    4  = 44428
    5  = 44459
    6  = 44489
    7  = 44498
    8  = 44520
    9  = 44561
    10 = 44581
    11 = 44613
    12 = 44824
    13 = 44854
    14 = 44885
    15 = 44926
}
foreach ($r in $cVals.Keys) {
    $ws.Range("C$r").Value = $cVals[$r]
}

# --- Append new rows 16-27 with the continued date series -------------------
$newCVals = @{
    16 = 44946
    17 = 44978
    18 = 45005
    19 = 45036
    20 = 45066
    21 = 45097
    22 = 45127
    23 = 45158
    24 = 45189
    25 = 45219
    26 = 45250
    27 = 45291
}

# Copy the formatting of the last existing data row down onto the new rows
# before filling them in, so the new cells pick up the same number formats /
# fills / alignment as the rest of the table.
$ws.Range("C15:G15").Copy()
$ws.Range("C16:G27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($r in ($newCVals.Keys | Sort-Object)) {
    $ws.Range("C$r").Value = $newCVals[$r]
}

# --- (Re)write the D/E/F/G formulas for every data row 4-27 -----------------
for ($r = 4; $r -le 27; $r++) {
    $ws.Range("D$r").Formula = "=(YEAR(`$B`$2)-YEAR(C$r))*12"
    $ws.Range("E$r").Formula = "=MONTH(`$B`$2)-MONTH(C$r)"
    $ws.Range("F$r").Formula = "=CONCATENATE(`"M`",D$r+E$r)"
    $ws.Range("G$r").Formula = "=CONCATENATE(`"M`",REPT(0,2-LEN(D$r+E$r))&(D$r+E$r))"
}

# --- Match the saved selection from the edit (cell I23) ---------------------
$ws.Range("I23").Select()
